$d = $word.ActiveDocument

# 1. Extend the text in the "(easy to code up ...)" paragraph.
#    Replace the tail portion "t filtering data before running)" with the
#    longer version that includes the new sentence about not doing new things.
$rngReplace = $d.Content
$rngReplace.Find.Execute("t filtering data before running)", $true, $false, $false, $false, $false, $true, 1, $false, "t filtering data before running. Don" + [char]0x2019 + "t do new things with finding similarity etc.)", 2)

# 2. Split that paragraph's single run into two runs: "(easy to code up – jus"
#    and "t filtering data before running. Don't do new things with finding
#    similarity etc.)" -- matching the target diff. We force the run boundary
#    by briefly adding (and then removing) a bookmark exactly at the split
#    point; Word (and this engine) splits runs around bookmark anchors, and
#    the split survives the bookmark's removal.
$splitRange = $d.Range(2175, 2175)
$d.Bookmarks.Add("TempSplitMark", $splitRange)
$d.Bookmarks.Item("TempSplitMark").Delete()

# 3. Move the "_GoBack" bookmark (which marks the location of the most
#    recent edit) from its old location -- right after "given good?" near
#    the end of the document -- to right after the text we just edited,
#    i.e. right after the closing parenthesis in "...etc.)".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$goBackRange = $d.Range(2257, 2257)
$d.Bookmarks.Add("_GoBack", $goBackRange)
